# Updates the cryptos price/volume table (GitHub Actions daily refresh).
# Cells whose new text is a plain number (e.g. "22.40", "0.0320") are written
# with a leading apostrophe so Excel keeps them as Text instead of silently
# coercing to a Number and dropping the trailing/insignificant zero, which
# would otherwise corrupt the string (matches how the source stores every
# Price/Volume cell as inline text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.043.62'
$ws.Range('E2').Value = '  +3.83%  '
$ws.Range('D3').Value = '2.656.20'
$ws.Range('E3').Value = '  +6.19%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = "'114.28"
$ws.Range('E5').Value = '  +7.65%  '
$ws.Range('D6').Value = "'326.51"
$ws.Range('E6').Value = '  +2.54%  '
$ws.Range('E7').Value = '  +2.00%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  +3.51%  '
$ws.Range('D10').Value = "'41.28"
$ws.Range('D11').Value = "'20.11"
$ws.Range('E11').Value = '  -0.88%  '
$ws.Range('D12').Value = "'0.0827"
$ws.Range('E12').Value = '  +2.90%  '
$ws.Range('E13').Value = '  +0.55%  '
$ws.Range('E14').Value = '  +3.97%  '
$ws.Range('D15').Value = '3.070.51'
$ws.Range('E15').Value = '  +6.26%  '
$ws.Range('D16').Value = '2.688.00'
$ws.Range('E16').Value = '  +7.94%  '
$ws.Range('D17').Value = "'0.877"
$ws.Range('E17').Value = '  +5.60%  '
$ws.Range('D18').Value = '49.983.71'
$ws.Range('E18').Value = '  +4.02%  '
$ws.Range('E19').Value = '  +3.17%  '
$ws.Range('D20').Value = "'6.79"
$ws.Range('E20').Value = '  +2.96%  '
$ws.Range('D21').Value = "'2.94"
$ws.Range('E21').Value = '  -1.93%  '
$ws.Range('E22').Value = '  +3.15%  '
$ws.Range('D23').Value = "'72.59"
$ws.Range('E23').Value = '  +2.00%  '
$ws.Range('D24').Value = "'277.02"
$ws.Range('E24').Value = '  +3.16%  '
$ws.Range('E25').Value = '  +2.97%  '
$ws.Range('D26').Value = "'26.97"
$ws.Range('E26').Value = '  +4.53%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').Value = "'10.05"
$ws.Range('E28').Value = '  +2.89%  '
$ws.Range('D29').Value = "'36.97"
$ws.Range('E29').Value = '  +6.49%  '
$ws.Range('E30').Value = '  +1.22%  '
$ws.Range('E31').Value = '  +1.42%  '
$ws.Range('D32').Value = "'50.17"
$ws.Range('E32').Value = '  +1.64%  '
$ws.Range('E33').Value = '  +3.85%  '
$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D34').Value = "'19.51"
$ws.Range('E34').Value = '  +1.79%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = "'0.0817"
$ws.Range('E35').Value = '  +5.60%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').Value = "'5.02"
$ws.Range('E37').Value = '  +8.97%  '
$ws.Range('E38').Value = '  +6.50%  '
$ws.Range('D40').Value = "'124.63"
$ws.Range('E40').Value = '  +1.08%  '
$ws.Range('E41').Value = '  +2.21%  '
$ws.Range('D42').Value = "'22.40"
$ws.Range('E42').Value = '  +0.65%  '
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('D44').Value = "'0.0320"
$ws.Range('E44').Value = '  +5.53%  '
$ws.Range('D45').Value = '2.099.42'
$ws.Range('E45').Value = '  +4.80%  '
$ws.Range('E46').Value = '  +5.84%  '
$ws.Range('E47').Value = '  +13.10%  '
$ws.Range('E48').Value = '  +4.12%  '
$ws.Range('E49').Value = '  +2.14%  '
$ws.Range('D50').Value = "'5.37"
$ws.Range('E50').Value = '  +3.09%  '
$ws.Range('D51').Value = "'60.49"
$ws.Range('E51').Value = '  +6.35%  '
